# Data updated in excel sheet
#
# - AccountCreationData!A2:A4 hold test email addresses that get
#   incremented (newtest1/2/3 -> newtest4/5/6@gmail.com).
# - The AccountCreationData sheet becomes the active sheet/tab, with the
#   selection moved to B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountCreationData")

$ws.Range("A2").Value = "newtest4@gmail.com"
$ws.Range("A3").Value = "newtest5@gmail.com"
$ws.Range("A4").Value = "newtest6@gmail.com"

# Make this sheet the active one (moves tabSelected + workbook activeTab)
# and move the cell selection/cursor to B13.
$ws.Activate()
$ws.Range("B13").Select()
